$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (C) column dates for all data rows (2-9) from 46073 -> 46074
foreach ($r in 2..9) {
    $ws.Cells.Item($r, 3).Value2 = 46074
}

# Swap rows 3 and 4 (Beteckning/Datum/Area values)
$a3 = $ws.Cells.Item(3, 1).Value2
$b3 = $ws.Cells.Item(3, 2).Value2
$g3 = $ws.Cells.Item(3, 7).Value2

$a4 = $ws.Cells.Item(4, 1).Value2
$b4 = $ws.Cells.Item(4, 2).Value2
$g4 = $ws.Cells.Item(4, 7).Value2

$ws.Cells.Item(3, 1).Value2 = $a4
$ws.Cells.Item(3, 2).Value2 = $b4
$ws.Cells.Item(3, 7).Value2 = $g4

$ws.Cells.Item(4, 1).Value2 = $a3
$ws.Cells.Item(4, 2).Value2 = $b3
$ws.Cells.Item(4, 7).Value2 = $g3

# Swap rows 5 and 6 (Beteckning/Datum/Area values)
$a5 = $ws.Cells.Item(5, 1).Value2
$b5 = $ws.Cells.Item(5, 2).Value2
$g5 = $ws.Cells.Item(5, 7).Value2

$a6 = $ws.Cells.Item(6, 1).Value2
$b6 = $ws.Cells.Item(6, 2).Value2
$g6 = $ws.Cells.Item(6, 7).Value2

$ws.Cells.Item(5, 1).Value2 = $a6
$ws.Cells.Item(5, 2).Value2 = $b6
$ws.Cells.Item(5, 7).Value2 = $g6

$ws.Cells.Item(6, 1).Value2 = $a5
$ws.Cells.Item(6, 2).Value2 = $b5
$ws.Cells.Item(6, 7).Value2 = $g5
